$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the old French-only row ("Les liens vers des fichiers telechargeables...",
# row 31) so every row below it shifts up by one.
$ws.Range("A31").EntireRow.Delete()

# --- Replace the French criteria text with the new English wording ---
$ws.Range("A1").Value = "Criteria"
$ws.Range("B1").Value = "Check? Y / N"
$ws.Range("C1").Value = "Comments"

$ws.Range("A2").Value = "Slide Layout"
$ws.Range("A3").Value = "Use the Slide Master view: select the View tab > Slide Master to create or modify slides layouts"
$ws.Range("A4").Value = "Create slide layouts for every type of slide"

$ws.Range("A5").Value = "Properties"
$ws.Range("A6").Value = "Title, author and tags fields have to be completed in tab File > Info > Properties."
$ws.Range("A7").Value = "The main language must be defined in the tab File > Options > Language"
$ws.Range("A8").Value = "Every punctual change of language must be identified by selecting the text and: tab Review > Language > Set proofing language"

$ws.Range("A9").Value = "Titles"
$ws.Range("A10").Value = "Give every slide a unique and descriptive title"

$ws.Range("A11").Value = "Wording"
$ws.Range("A12").Value = "If necessary, keep accentuation on capital letters (example: É or Ç): use the Insert tab > Symbol"
$ws.Range("A13").Value = "Avoid abbreviations and acronyms; explain them at the first occurrence"

$ws.Range("A14").Value = "Formatting / layout"
$ws.Range("A15").Value = "Use familiar sans serif fonts (such as Arial or Calibri), in a larger font size (18pt or larger)."
$ws.Range("A16").Value = "Begin all sentences by a capital letter but avoid using all capital letters and excessive italics or underlines"
$ws.Range("A17").Value = "Avoid justifying the text; prefer left alignment."
$ws.Range("A18").Value = "In order to add white space between sentences and paragraphs, use option: right click on the text and select Paragraph… > Indents and Spacing"

$ws.Range("A19").Value = "Colors and contrasts"
$ws.Range("A20").Value = "The presentation of text should have a contrast ratio of at least 4.5:1`nyou can download the Colour Contrast Analyser tool"
$ws.Range("D20").Value = "Colour Contrast Analyser"
$ws.Range("A21").Value = "Ensure that color is not the only means of conveying information"

$ws.Range("A22").Value = "Bulleted Lists or Numbered Lists"
$ws.Range("A23").Value = "the appropriate and native functions are used in order to create lists"

$ws.Range("A24").Value = "Pictures"
$ws.Range("A25").Value = "alt text: it conveys the content and the purpose of the image and its context in a concise and unambiguous manner."
$ws.Range("A26").Value = "If the picture is only decorative, leave the Description text box blank"
$ws.Range("A27").Value = "graphics / diagrams:`n- SmartArt instead of images (Insert> SmartArt) and alt text is completed"

$ws.Range("A28").Value = "Links"
$ws.Range("A29").Value = "Hyperlinks titles are simples and concises"
$ws.Range("A30").Value = "For every downloadable file, it is necessary to specify its name, type, size and language (if different of the presentation language)"

$ws.Range("A31").Value = "Tables"
$ws.Range("A32").Value = "Tables are not used for layout purposes"
$ws.Range("A33").Value = "Ensure that tables don't contain split cells, merged cells, nested tables, or completely blank rows or columns."
$ws.Range("A34").Value = "If a table is too complex, you will need to provide an additional description: Right-click Format Shape > Alt Text."

$ws.Range("A35").Value = "Audio or video file"
$ws.Range("A36").Value = "Audio content must have a transcript"
$ws.Range("A37").Value = "Video content must have synchronized captions and also needs a rich description of the visual components."
$ws.Range("A38").Value = "Avoid content that blinks or flashes more than three times a second."
$ws.Range("A39").Value = "Avoid automatic launching or slide transitions"

$ws.Range("A40").Value = "Final Checks"
$ws.Range("A41").Value = "run Accessibility Checker to make sure your content is accessible.`nIn File tab > Check for Issues (copy the result in the corresponding tab)    "
$ws.Range("A42").Value = "Check the order in which the screen readers read the slide contents"

$ws.Range("A43").Value = "Convert to PDF"
$ws.Range("A44").Value = "Save you presentation as PDF: File tab > Save As and save as type PDF. Check the « Document structure tags for accessibility » option."
$ws.Range("A45").Value = "Note: To check PDF document accessibility: download PDF Accessibility Checker (PAC 3)."
$ws.Range("A46").Value = "Note: If your presentation have more than 50 slides, it is better to advise your readers to change the following setting in Adobe: Edit > Preferences > Reading > Screen Reader Options > Page vs Document: `"Read the entire document`""

# --- Row-height tweaks that come with the new (longer/shorter) English wording ---
$ws.Rows.Item(3).RowHeight = 31.5
$ws.Rows.Item(8).RowHeight = 31.5
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 31.5
$ws.Rows.Item(13).RowHeight = 15.75
$ws.Rows.Item(15).RowHeight = 31.5
$ws.Rows.Item(20).RowHeight = 31.5
$ws.Rows.Item(30).RowHeight = 31.5

# --- Update the view: scrolled down a bit further, selection now on the Properties header ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A5:C5").Select()
